$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C (Response) to match the new content being added.
$ws.Columns.Item(3).ColumnWidth = 468.17

# Append four new Q&A rows (60-64) beneath the existing data.
$ws.Range("A60").Value = "What are the number of curves that can be plotted in GEO?"
$ws.Range("B60").Value = "llama3.2:latest"
$ws.Range("C60").Value = "According to the information provided, there is no specific limit on the number of curves that can be plotted in GEO. However, it does mention that ""Number of curve shades per plot"" has a value of 250. This suggests that while you may not be able to display all possible curves at once, there is no hard limit on the total number of curves available for plotting."

$ws.Range("A61").Value = "How many curve shades can I create?"
$ws.Range("B61").Value = "llama3.2:latest"
$ws.Range("C61").Value = "According to the document, you can create 250 curve shades."

$ws.Range("A62").Value = "What are the types of limits in GEO data?"
$ws.Range("B62").Value = "llama3.2:latest"
$ws.Range("C62").Value = "In GEO, there are two main types of limits: `n1. Modifiers`n2. Lithologies`n3. Symbols`n4. Texts`n5. Lines"

$ws.Range("A63").Value = "What is an Output Database Type file?"
$ws.Range("B63").Value = "llama3.2:latest"
$ws.Range("C63").Value = "The Output Database File, also known as ODF (Output Database File), is a compact database file created by GEO that contains all the information needed for rapid recall of data related to a particular well. It includes loaded curve data, layout-related information, geological interpretation data, user-inputted information, libraries and fonts, embedded objects, and links to external documents."

$ws.Range("A64").Value = "Can you export image files?"
$ws.Range("B64").Value = "llama3.2:latest"
$ws.Range("C64").Value = "Yes, you can export image files from GEO. To do this:`n1. Go to the Exporttab.`n2. Select HTML Using... and choose your preferred image format (e.g., EMF).`n3. An ExportTodialog box will open where you can enter a Filename and select the directory where it will be stored.`nNote that exporting to HTML outputs an image of your log using spliced web-supported images, which are outputted to a location specified by the user at the point of saving, keeping them together."

# Restore default row heights so rows whose multi-line text made the host
# auto-grow them stay at the sheet's normal (non-custom) height.
$ws.Rows.Item(60).AutoFit()
$ws.Rows.Item(61).AutoFit()
$ws.Rows.Item(62).AutoFit()
$ws.Rows.Item(63).AutoFit()
$ws.Rows.Item(64).AutoFit()
